# Eagle-related works slide: split the "Privacy Related" ellipse label
# into two separate paragraphs: "Privacy" and "Sensitive".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(10)
$tr = $shp.TextFrame.TextRange

$tr.Text = "Privacy"
$tr.InsertAfter([char]13 + "Sensitive")
